# Apply crypto price/volume updates to the worksheet, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.366.91'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.63%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.008.65'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.45%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.13'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.06%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.53'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.71%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.19%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.63%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.002.91'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.47%  '
# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.35%  '
# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +8.24%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.52%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.76'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.28%  '
# Row 15
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.29%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.509.11'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.63%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.32'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +7.63%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.000.87'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.25%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '59.506.18'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.88%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '431.90'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.08%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.75'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.91%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.723'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.35%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.14'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.27%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.44'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.22%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.80'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.55%  '
# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.04%  '
# Row 27
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.21'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +10.40%  '
# Row 28
$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.995'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.43%  '
# Row 29
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.07%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.88'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.91%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.84'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.18%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.10'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.49%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0999'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.03%  '
# Row 34
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.98'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.65%  '
# Row 35
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.29%  '
# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +10.65%  '
# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.70%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.88'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.74%  '
# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.54%  '
# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.15%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '404.26'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.97%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0354'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.11%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.767.34'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.19%  '
# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.97%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.253'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.38%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '35.90'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +26.13%  '
# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.03%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.42'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.89%  '
# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.01%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.02'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.71%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.58'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.16%  '

Write-Output "Applied 95 cell updates to sheet1."
